# Insert a new weekly data row at row 27, pushing existing rows 27:56 down to
# 28:57 (mirrors the "new row added" shape of the diff: old row 27 -> new row
# 28, ..., old row 56 -> new row 57), then populate the new row 27 with the
# latest week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything at/after row 27 down by one row.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 (same static attributes as every other
# row in this data set; only date / volume / price columns differ).
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44533
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112038
$ws.Cells.Item(27, 7).Value = "Cebollín baby"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 250
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 1100
$ws.Cells.Item(27, 13).Value = 1050
$ws.Cells.Item(27, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 525
$ws.Cells.Item(27, 17).Value = 2
$ws.Cells.Item(27, 18).Value = "Hortaliza"
